# "Logged Week 16 and performed season sim from Week 17"
#
# Rushing sheet: a handful of weekly totals (weeks already logged) get
# bumped up. Receiving sheet: a new Week-16 row is logged for A.Brown and
# every later row shifts down one (the old trailing A.Brown aggregate row
# is gone, replaced by this newly logged week), plus a handful of other
# weekly totals change from the re-simmed remainder of the season.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Rushing" - update the handful of changed weekly totals
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

$rushing.Range("C2").Value = 7
$rushing.Range("D2").Value = 6
$rushing.Range("E2").Value = 10

$rushing.Range("C3").Value = 13
$rushing.Range("D3").Value = 8
$rushing.Range("E3").Value = 6
$rushing.Range("F3").Value = 5

$rushing.Range("C6").Value = 54
$rushing.Range("D6").Value = 36
$rushing.Range("F6").Value = 22

$rushing.Range("C7").Value = 30
$rushing.Range("D7").Value = 20

$rushing.Range("C27").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Receiving" - rewrite the full player log: a new Week 16 row for
# A.Brown is inserted (row 7), pushing every following player's row down
# by one, and the old trailing A.Brown summary row (old row 17) is
# dropped since it's superseded by the newly logged week.
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

function Set-Row($ws, $row, $week, $name, $c, $d, $e, $f, $g, $h) {
    $ws.Range("A$row").Value = $week
    $ws.Range("B$row").Value = $name
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

Set-Row $receiving 2  0  "J.McNichols"        15 9  1  1  5  3
Set-Row $receiving 3  1  "K.Blasingame"       1  1  0  0  0  0
Set-Row $receiving 4  2  "D.Evans"            2  2  0  0  0  0
Set-Row $receiving 5  3  "D.Foreman"          8  7  1  1  1  1
Set-Row $receiving 6  4  "D.Hilliard"         22 15 0  0  5  4
Set-Row $receiving 7  5  "A.Brown"            93 70 22 16 11 9
Set-Row $receiving 8  7  "C.Rogers"           27 18 4  2  2  0
Set-Row $receiving 9  8  "R.McMath"           2  2  1  0  1  1
Set-Row $receiving 10 9  "M.Johnson"          13 6  6  3  3  1
Set-Row $receiving 11 10 "D.Fitzpatrick"      6  5  2  0  2  2
Set-Row $receiving 12 11 "N.Westbrook-Ikhine" 24 19 7  3  4  4
Set-Row $receiving 13 12 "C.Hollister"        6  4  1  0  2  1
Set-Row $receiving 14 13 "A.Firkser"          29 24 2  0  2  1
Set-Row $receiving 15 14 "M.Pruitt"           14 10 3  2  6  3
Set-Row $receiving 16 15 "G.Swaim"            29 23 1  1  6  3
Set-Row $receiving 17 16 "T.Hudson"           2  1  1  0  0  0

# Row 17 ("A16" week-index cell) used to carry a special thin-border-only
# style (distinct from the boxed style every other week-index cell in
# column A uses). That row is now just a normal logged week, so drop the
# special style by copying the uniform format from row 16's A cell.
$receiving.Range("A16").Copy() | Out-Null
$receiving.Range("A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$receiving.Range("J10").Select() | Out-Null
